$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Mobile No" column (F) originally stored each phone number as *text*
# (via a shared string). The sheet was edited so that column F now holds a
# (slightly re-typed/truncated) *numeric* phone number, and the original,
# full phone number is preserved as a number in a new column I.
# ---------------------------------------------------------------------------

# Row 14 - Devika BR
$ws.Range("F14").Value = 940030761
$ws.Range("I14").Value = 9400301761

# Row 16 - Gowri mohan
$ws.Range("F16").Value = 894362081
$ws.Range("I16").Value = 8943862081

# Row 27 - Gouri M Biju
$ws.Range("F27").Value = 776702081
$ws.Range("I27").Value = 7736702081

# Row 38 - Archana S
$ws.Range("F38").Value = 940930260
$ws.Range("I38").Value = 9400930260

# Row 50 - Aswathy Ashok
$ws.Range("F50").Value = 812957440
$ws.Range("I50").Value = 8129574840

# Row 52 - Nesla
$ws.Range("F52").Value = 859092902
$ws.Range("I52").Value = 8590929802

# Row 55 - Abhiram m p
$ws.Range("F55").Value = 751027418
$ws.Range("I55").Value = 7510274318

# Row 61 - Abhiram Renjith (text -> number, no extra column here)
$ws.Range("F61").Value = 8089311685

# Row 75 - Arun Krishna K U
$ws.Range("F75").Value = 963082570
$ws.Range("I75").Value = 9633082570

# Row 79 - Ebrahim Roshan U
$ws.Range("F79").Value = 996184053
$ws.Range("I79").Value = 9961840533

# ---------------------------------------------------------------------------
# Rows 67:72 of the running-id column H had individual, non-shared formulas
# (=H66+1, =H67+1, ...). Re-entering the same formula across the whole
# range collapses them back into a single shared-formula group, matching
# how the rest of column H already behaves.
# ---------------------------------------------------------------------------
$ws.Range("H67:H72").Formula = "=H66+1"

# ---------------------------------------------------------------------------
# Drop the stale column-outline level that no longer groups anything.
# ---------------------------------------------------------------------------
$ws.Columns("A:A").Ungroup()

# ---------------------------------------------------------------------------
# Update the active selection to reflect where the editor left off.
# ---------------------------------------------------------------------------
$ws.Range("I22").Select()
